$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.547.72"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "1.755.67"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'324.50"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("E7").Value = "  +1.68%  "

$ws.Range("E8").Value = "  -1.82%  "

$ws.Range("D9").Value = "'0.07468"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("D10").Value = "'41.52"
$ws.Range("E10").Value = "  -1.43%  "

$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").Value = "'20.82"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("E14").Value = "  -0.68%  "

$ws.Range("D15").Value = "'7.150"
$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("D16").Value = "1.752.60"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").Value = "'93.63"
$ws.Range("E17").Value = "  +0.71%  "

$ws.Range("D18").Value = "'0.00001053"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").Value = "'0.06411"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D22").Value = "'5.737"
$ws.Range("E22").Value = "  -2.01%  "

$ws.Range("D23").Value = "27.592.93"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").Value = "'11.23"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("D26").Value = "'165.84"
$ws.Range("E26").Value = "  +2.47%  "

$ws.Range("D27").Value = "'20.13"
$ws.Range("E27").Value = "  -1.20%  "

$ws.Range("D28").Value = "1.953.32"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").Value = "'2.133"
$ws.Range("E29").Value = "  +0.47%  "

$ws.Range("D30").Value = "'125.46"
$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("D31").Value = "'1.080"
$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "'5.526"
$ws.Range("E34").Value = "  -0.56%  "

$ws.Range("E35").Value = "  -2.27%  "

$ws.Range("D36").Value = "'0.02283"
$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("D39").Value = "'0.6277"
$ws.Range("E39").Value = "  -1.74%  "

$ws.Range("D40").Value = "'4.921"
$ws.Range("E40").Value = "  -1.40%  "

$ws.Range("D41").Value = "'1.182"
$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("D42").Value = "'1.386"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").Value = "'7.786"
$ws.Range("E43").Value = "  -0.22%  "

$ws.Range("D44").Value = "'13.17"
$ws.Range("E44").Value = "  -0.78%  "

$ws.Range("D46").Value = "'0.5865"
$ws.Range("E46").Value = "  -0.76%  "

$ws.Range("D47").Value = "'122.09"
$ws.Range("E47").Value = "  +0.26%  "

$ws.Range("D48").Value = "'1.937"
$ws.Range("E48").Value = "  -1.20%  "

$ws.Range("D49").Value = "'0.06899"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("E50").Value = "  -2.73%  "

$ws.Range("D51").Value = "'72.25"
$ws.Range("E51").Value = "  -0.39%  "
